# Apply the NCM results re-ordering (newest created-date first) + new rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D to fit the longer descriptions now present in the data.
$ws.Columns.Item(4).ColumnWidth = 129

# Column C holds NCM codes that look numeric (e.g. "0102.3", "01.05") - keep as text
# (data rows only, the header "NCM" in C1 is untouched).
$ws.Range("C2:C9").NumberFormat = "@"

# Full desired contents for data rows 2-9 (header in row 1 stays as-is).
$rows = @(
    @{ A = "0adac61f-b078-4812-86d4-af98abc9538f"; B = ""  ; C = "0102.3"     ; D = "BÚFALOS:"                                                                                                                           ; E = "2025-05-21 10:56:42.301000" },
    @{ A = "8ed7eaf9-259f-4422-8a86-6bf706a844c6"; B = "NT"; C = "0102.29.11" ; D = "Prenhes ou com cria ao pé"                                                                                                          ; E = "2025-05-21 10:56:42.251000" },
    @{ A = "e1264496-3d8a-4650-8ee2-829419786d9e"; B = ""  ; C = "01.05"      ; D = "AVES DA ESPÉCIE GALLUS DOMESTICUS, PATOS, GANSOS, PERUS, PERUAS E GALINHAS-D’ANGOLA (PINTADAS), DAS ESPÉCIES DOMÉSTICAS, VIVOS."   ; E = "2025-05-20 19:06:22.392000" },
    @{ A = "3831aba8-4ff2-4c7c-859c-96789b1675c7"; B = ""  ; C = "01.05"      ; D = "AVES DA ESPÉCIE GALLUS DOMESTICUS, PATOS, GANSOS, PERUS, PERUAS E GALINHAS-D’ANGOLA (PINTADAS), DAS ESPÉCIES DOMÉSTICAS, VIVOS."   ; E = "2025-05-20 18:57:28.911000" },
    @{ A = "f768e507-e1c8-4764-ad06-cfc6ff08999a"; B = ""  ; C = "0102.39.1"  ; D = "PARA REPRODUÇÃO"                                                                                                                    ; E = "2025-05-20 18:57:28.495000" },
    @{ A = "105f25df-4ffc-4487-b3bb-aae17a07d567"; B = ""  ; C = "0102.39.1"  ; D = "PARA REPRODUÇÃO"                                                                                                                    ; E = "2025-05-20 18:53:10.311000" },
    @{ A = "e58092ef-e69a-47ec-9b6c-67245eac6716"; B = ""  ; C = "1"          ; D = "ANIMAIS VIVOS"                                                                                                                      ; E = "2025-05-20 18:41:05.422000" },
    @{ A = "c1d3f1f1-0422-4883-a063-75d18b4d26d1"; B = ""  ; C = "1"          ; D = "ANIMAIS VIVOS"                                                                                                                      ; E = "2025-05-20 16:16:30.735000" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $r++
}
